# Generate Report for Archive
#
# Localization status moved from "Ready for handoff" to "In Translation"
# across the Overview summary sheet (columns E/F, one per locale) and the
# corresponding per-locale detail sheets (column C, "Status"). The now
# shorter status text no longer needs as wide a column, so the related
# Status columns are narrowed to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 12.5   # renders as the narrower column width used for Status columns

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn detail sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de detail sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
